# Update "想去人数" (want-to-go count) figures in column F
# for the "展览" (Exhibitions) and "全部类型" (All Types) sheets.
# Both sheets carry the same data table, so the same F-column
# updates are applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2 = 1720
    4 = 481
    6 = 80
    7 = 650
    8 = 411
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
